$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), reusing the same style as
# the existing header cells (e.g. H1: bold, bordered, centered).
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J, rows 2-10 (plain numbers, no special style).
$iValues = @(7, 8, 7, 7, 8, 9, 4, 7, 9)
$jValues = @(7, 8, 7, 8, 8, 9, 4, 7, 9)

for ($r = 0; $r -lt 9; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
